$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 78
$prevRow = 77

# Copy the formatting (styles) of the previous data row into the new row
$ws.Range("A$prevRow`:V$prevRow").Copy()
$ws.Range("A$row`:V$row").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 77
$ws.Cells.Item($row, 2).Value = "poland"
$ws.Cells.Item($row, 3).Value = "ekstraklasa"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45196.83333333334
$ws.Cells.Item($row, 6).Value = "Pogon Szczecin"
$ws.Cells.Item($row, 7).Value = 3
$ws.Cells.Item($row, 8).Value = "Legia"
$ws.Cells.Item($row, 9).Value = 4
$ws.Cells.Item($row, 10).Value = 2.37
$ws.Cells.Item($row, 11).Value = "24/09/2023 16:42"
$ws.Cells.Item($row, 12).Value = 2.41
$ws.Cells.Item($row, 13).Value = "27/09/2023 19:50"
$ws.Cells.Item($row, 14).Value = 3.5
$ws.Cells.Item($row, 15).Value = "24/09/2023 16:42"
$ws.Cells.Item($row, 16).Value = 3.75
$ws.Cells.Item($row, 17).Value = "27/09/2023 19:25"
$ws.Cells.Item($row, 18).Value = 2.83
$ws.Cells.Item($row, 19).Value = "24/09/2023 16:42"
$ws.Cells.Item($row, 20).Value = 2.84
$ws.Cells.Item($row, 21).Value = "27/09/2023 19:50"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/pogon-szczecin-legia/rypDtbZ5/"
